$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.286.58'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.618.06'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.48'
$ws.Range('E5').Value = '  +5.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.29'
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  +5.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.83'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  +5.67%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.076.52'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '59.212.65'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.21'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.622.20'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.47'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.59'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.18'
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.20'
$ws.Range('E21').Value = '  +0.41%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.19'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.430'
$ws.Range('E24').Value = '  +4.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.166'
$ws.Range('E25').Value = '  +2.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.22'
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0769'
$ws.Range('E28').Value = '  +5.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('E31').Value = '  +4.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.91'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.96'
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('E34').Value = '  +2.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.915'
$ws.Range('E35').Value = '  +12.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.905'
$ws.Range('E36').Value = '  +9.72%  '
$ws.Range('E37').Value = '  +2.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.24'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.48'
$ws.Range('E39').Value = '  +3.84%  '
$ws.Range('E40').Value = '  +1.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '286.51'
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.996'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.603'
$ws.Range('E43').Value = '  +2.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0543'
$ws.Range('E44').Value = '  +2.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0960'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.72'
$ws.Range('E47').Value = '  +5.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0229'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.957.28'
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '117.86'
$ws.Range('E50').Value = '  +5.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.18'
$ws.Range('E51').Value = '  +1.54%  '
